$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (date moved from 2022-06-21 to 2022-06-22)
$ws.Name = "Through 2022-06-22"

# Update the "June" row label text
$ws.Range("A7").Value = "June (through 06-22)"

# Update June row (row 7) values
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 25
$ws.Range("D7").Value = 54
$ws.Range("E7").Value = 43
$ws.Range("F7").Value = 34
$ws.Range("G7").Value = 84
$ws.Range("H7").Value = 88
$ws.Range("I7").Value = 104

# Update Total row (row 8) values
$ws.Range("B8").Value = 120
$ws.Range("C8").Value = 234
$ws.Range("D8").Value = 370
$ws.Range("E8").Value = 338
$ws.Range("F8").Value = 238
$ws.Range("G8").Value = 442
$ws.Range("H8").Value = 719
$ws.Range("I8").Value = 767
